# Minor changes - use case #1
# Applies the edits described in the commit diff to the Swimclub use-case
# document: tidy up the Scope text/bookmark, italicise a phrase in the
# main flow, and renumber/reword several "Extensions" bullet points.

$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs.Item($ParaIndex)
    $searchRange = $d.Range($p.Range.Start, $p.Range.End)
    $found = $searchRange.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    return $found
}


# --- 1. Scope row: "Delfinklub management system" -----------------------
# Collapses " management" + bookmark + " system" into a single run,
# dropping the now-redundant _GoBack bookmark from this cell.
$pScope = $d.Paragraphs.Item(8)
$scopeRange = $d.Range($pScope.Range.Start, $pScope.Range.End)
$scopeRange.Text = "Delfinklub management system"

# --- 2. Main flow bullet 1: italicise "opret nyt medlem" ------------------
$pFlow1 = $d.Paragraphs.Item(26)
$flow1Full = $d.Range($pFlow1.Range.Start, $pFlow1.Range.End)
$flow1Full.Text = "Formanden vælger opret nyt medlem fra en menu."
$pFlow1b = $d.Paragraphs.Item(26)
$italicRange = $d.Range($pFlow1b.Range.Start, $pFlow1b.Range.End)
$italicRange.Find.Execute("opret nyt medlem")
$italicRange.Font.Italic = $true

# --- 3. Extensions renumbering -------------------------------------------
# "4a. System prompter for gyldigt navn." -> "3b. System prompter for gyldigt navn."
Replace-InParagraph 38 "4a." "3b."

# "5b. Formand indtaster ikke en numerisk værdi." -> "5a. Formand indtaster ikke en numerisk værdi."
Replace-InParagraph 40 "5b." "5a."

# "6b. System prompter for gyldig alder." -> "5b. System prompter for gyldig alder."
Replace-InParagraph 41 "6b." "5b."

# "7c. Formand vælger en ikkeeksisterende aktivitetsform." ->
# "7a. Formand vælger en ikkeeksisterende mulighed."
Replace-InParagraph 43 "7c." "7a."
Replace-InParagraph 43 "ikkeeksisterende aktivitetsform" "ikkeeksisterende mulighed"

# "8c. System prompter for eksisterende aktivitetsform" ->
# "7b. System prompter for eksisterende muligheder."
Replace-InParagraph 44 "8c." "7b."
Replace-InParagraph 44 "eksisterende aktivitetsform" "eksisterende muligheder."

# Italic note: "...eksisterende aktivitetsform er valgt." -> "...eksisterende mulighed er valgt."
Replace-InParagraph 45 "eksisterende aktivitetsform" "eksisterende mulighed"
